$wb = $excel.ActiveWorkbook

# New record to append (Kadastro App: Yeni kayit eklendi: 2991)
$newRow = @("2991", "2025-09-10", "Erdemli", "1", "ÇAP", "SEVİL SARAÇER (Tekniker)")

# Both the main "Kayitlar" list and the per-birim "Erdemli" list carry this row.
$sheetNames = @("Kayitlar", "Erdemli")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $targetRow = 57

    # Force the new row to be stored as text (matches every existing cell in
    # the sheet, which are all text values even when they look numeric),
    # so values like "2991" / "1" / "2025-09-10" aren't reinterpreted as a
    # number or a date.
    $rowRange = $ws.Range("A" + $targetRow + ":F" + $targetRow)
    $rowRange.NumberFormat = "@"

    $ws.Cells.Item($targetRow, 1).Value = $newRow[0]
    $ws.Cells.Item($targetRow, 2).Value = $newRow[1]
    $ws.Cells.Item($targetRow, 3).Value = $newRow[2]
    $ws.Cells.Item($targetRow, 4).Value = $newRow[3]
    $ws.Cells.Item($targetRow, 5).Value = $newRow[4]
    $ws.Cells.Item($targetRow, 6).Value = $newRow[5]
}
